# NatMI LR-pair export (Nppa -> Npr2) refreshed with new TPM numbers.
#
# The previous export had three "Sending cluster" groups (ECs, FAPs,
# Inflammatory-Mac) x three "Target cluster" values (ECs, FAPs, MuSCs) = 9
# data rows. The refreshed run drops the "ECs" sending-cluster group
# entirely and recomputes the metrics for the remaining two groups, so the
# sheet shrinks from 9 data rows (A1:T10) to 6 data rows (A1:T7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old rows 2-4 (Sending cluster = "ECs"); rows 5-10 shift up to
# become rows 2-7.
$ws.Rows("2:4").Delete()

# Updated data, keyed by the final row number after the delete above.
$rows = @(
    @{ Row=2; A="FAPs";             D="ECs";   E=1; F=0.3333333333333333; G=0.2077686666666667; H=0.623306;  I=0.3307450952508051; J=0.3307450952508051; M=17.160689;         N=51.482067;         O=0.5879310549011342; P=0.5879310549011342; Q=3.565453472611334;  R=32.089081253502;   S=0.1944553127541819;  T=0.1944553127541819 }
    @{ Row=3; A="FAPs";             D="FAPs";  E=1; F=0.3333333333333333; G=0.2077686666666667; H=0.623306;  I=0.3307450952508051; J=0.3307450952508051; M=8.435839666666666; N=25.307519;         O=0.289014742601545;  P=0.289014742601545;  Q=1.752703159757111;  R=15.774328437814;   S=0.09559020857063492; T=0.09559020857063492 }
    @{ Row=4; A="FAPs";             D="MuSCs"; E=1; F=0.3333333333333333; G=0.2077686666666667; H=0.623306;  I=0.3307450952508051; J=0.3307450952508051; M=3.591739;          N=10.775217;         O=0.1230542024973208; P=0.1230542024973208; Q=0.7462508230446667; R=6.716257407402;    S=0.04069957392598821; T=0.04069957392598821 }
    @{ Row=5; A="Inflammatory-Mac"; D="ECs";   E=2; F=0.6666666666666666; G=0.420415;           H=1.261245;  I=0.6692549047491948; J=0.6692549047491948; M=17.160689;         N=51.482067;         O=0.5879310549011342; P=0.5879310549011342; Q=7.214611065935;     R=64.931499593415;   S=0.3934757421469522;  T=0.3934757421469522 }
    @{ Row=6; A="Inflammatory-Mac"; D="FAPs";  E=2; F=0.6666666666666666; G=0.420415;           H=1.261245;  I=0.6692549047491948; J=0.6692549047491948; M=8.435839666666666; N=25.307519;         O=0.289014742601545;  P=0.289014742601545;  Q=3.546553533461666;  R=31.918981801155;   S=0.1934245340309101;  T=0.1934245340309101 }
    @{ Row=7; A="Inflammatory-Mac"; D="MuSCs"; E=2; F=0.6666666666666666; G=0.420415;           H=1.261245;  I=0.6692549047491948; J=0.6692549047491948; M=3.591739;          N=10.775217;         O=0.1230542024973208; P=0.1230542024973208; Q=1.510020951685;     R=13.590188565165;   S=0.08235462857133254; T=0.08235462857133254 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value  = $r.A        # Sending cluster
    $ws.Cells.Item($n, 2).Value  = "Nppa"      # Ligand symbol (unchanged)
    $ws.Cells.Item($n, 3).Value  = "Npr2"      # Receptor symbol (unchanged)
    $ws.Cells.Item($n, 4).Value  = $r.D        # Target cluster
    $ws.Cells.Item($n, 5).Value  = $r.E        # Ligand-expressing cells
    $ws.Cells.Item($n, 6).Value  = $r.F        # Ligand detection rate
    $ws.Cells.Item($n, 7).Value  = $r.G        # Ligand average expression value
    $ws.Cells.Item($n, 8).Value  = $r.H        # Ligand total expression value
    $ws.Cells.Item($n, 9).Value  = $r.I        # Ligand derived specificity (avg)
    $ws.Cells.Item($n, 10).Value = $r.J        # Ligand derived specificity (total)
    $ws.Cells.Item($n, 11).Value = 3           # Receptor-expressing cells (unchanged)
    $ws.Cells.Item($n, 12).Value = 1           # Receptor detection rate (unchanged)
    $ws.Cells.Item($n, 13).Value = $r.M        # Receptor average expression value
    $ws.Cells.Item($n, 14).Value = $r.N        # Receptor total expression value
    $ws.Cells.Item($n, 15).Value = $r.O        # Receptor derived specificity (avg)
    $ws.Cells.Item($n, 16).Value = $r.P        # Receptor derived specificity (total)
    $ws.Cells.Item($n, 17).Value = $r.Q        # Edge average expression weight
    $ws.Cells.Item($n, 18).Value = $r.R        # Edge total expression weight
    $ws.Cells.Item($n, 19).Value = $r.S        # Edge average expression derived specificity
    $ws.Cells.Item($n, 20).Value = $r.T        # Edge total expression derived specificity
}
